$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to remain a plain text string: Excel normally
    # auto-converts decimal-looking strings (e.g. "300.40") to numbers
    # when assigned via .Value, which would drop the trailing zero and
    # flip the cell to a Number type. Pre-formatting as text ("@")
    # keeps the literal text, then resetting the style back to Normal
    # avoids leaving a stray number-format style on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.658.57"
$ws.Range("E2").Value = "  +0.62%  "
Set-TextValue "D3" "2.285.51"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "300.40"
$ws.Range("E5").Value = "  -0.03%  "
Set-TextValue "D6" "99.01"
$ws.Range("E6").Value = "  +2.70%  "
Set-TextValue "D7" "0.499"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.511"
$ws.Range("E9").Value = "  +3.83%  "
Set-TextValue "D10" "35.71"
$ws.Range("E10").Value = "  +7.12%  "
Set-TextValue "D11" "0.0788"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  +2.10%  "
Set-TextValue "D13" "17.83"
$ws.Range("E13").Value = "  +11.43%  "
Set-TextValue "D14" "6.80"
$ws.Range("E14").Value = "  +1.34%  "
Set-TextValue "D15" "2.643.17"
$ws.Range("E15").Value = "  +0.31%  "
Set-TextValue "D16" "2.311.23"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  +0.83%  "
Set-TextValue "D18" "42.570.27"
$ws.Range("E18").Value = "  +0.60%  "
Set-TextValue "D19" "12.39"
$ws.Range("E19").Value = "  +5.56%  "
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("E21").Value = "  +0.33%  "
Set-TextValue "D22" "67.74"
$ws.Range("E22").Value = "  +1.86%  "
Set-TextValue "D23" "234.98"
$ws.Range("E23").Value = "  -0.42%  "
Set-TextValue "D24" "2.20"
$ws.Range("E24").Value = "  +11.91%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.77%  "
Set-TextValue "D27" "24.51"
$ws.Range("E27").Value = "  +2.96%  "
Set-TextValue "D28" "2.19"
$ws.Range("E28").Value = "  +0.92%  "
Set-TextValue "D29" "168.23"
$ws.Range("E29").Value = "  +0.54%  "
Set-TextValue "D30" "34.24"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -0.03%  "
Set-TextValue "D33" "4.96"
$ws.Range("E33").Value = "  +0.60%  "
Set-TextValue "D34" "17.45"
$ws.Range("E34").Value = "  +4.03%  "
Set-TextValue "D35" "4.57"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +1.85%  "
Set-TextValue "D40" "1.76"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("E41").Value = "  +0.11%  "
Set-TextValue "D42" "1.986.63"
$ws.Range("E42").Value = "  +1.40%  "
Set-TextValue "D43" "0.0287"
$ws.Range("E43").Value = "  +2.63%  "
Set-TextValue "D44" "2.22"
$ws.Range("E44").Value = "  -1.81%  "
Set-TextValue "D45" "10.08"
$ws.Range("E45").Value = "  +4.69%  "
Set-TextValue "D48" "55.37"
$ws.Range("E48").Value = "  +5.50%  "
Set-TextValue "D49" "2.512.03"
$ws.Range("E49").Value = "  +0.34%  "
Set-TextValue "D50" "1.52"
$ws.Range("E50").Value = "  +2.45%  "
Set-TextValue "D51" "4.48"
$ws.Range("E51").Value = "  -0.74%  "

# Rows 46/47: NEARProtocol moves up to rank 46, EnergySwap drops to rank 47
# (ranking/link/price/volume all change together for this pair of rows).
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D46" "2.87"
$ws.Range("E46").Value = "  +2.48%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "17.42"
$ws.Range("E47").Value = "  -1.14%  "

Write-Output "Applied cryptos update"
